$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("accountCreatonTest")
$ws.Activate()

# "Login Account creation page enhancements": the duplicate-data Student IDs
# used on the account-creation test sheet were renumbered (040216 -> 040516),
# and the two near-duplicate IDs used for row 2 (...aa) and row 3 (...ab)
# were consolidated into a single STD_040516aa value.

# Row 2 (QA-AC_CREATE-001): STD_040216aa -> STD_040516aa
$ws.Range("G2").Value = "STD_040516aa"
$ws.Range("J2").Value = "STD_040516aa"
$ws.Range("K2").Value = "STD_040516aa"
$ws.Range("L2").Value = "STD_040516aa"

# Row 3 (QA-AC_CREATE-002): STD_040216aa -> STD_040516aa, STD_040216ab -> STD_040516aa
$ws.Range("G3").Value = "STD_040516aa"
$ws.Range("J3").Value = "STD_040516aa"
$ws.Range("K3").Value = "STD_040516aa"
$ws.Range("L3").Value = "STD_040516aa"

# Row 4 (QA-AC_CREATE-003): STD_040216bb -> STD_040516ab, STD_040216aa -> STD_040516aa
$ws.Range("G4").Value = "STD_040516ab"
$ws.Range("J4").Value = "STD_040516aa"
$ws.Range("K4").Value = "STD_040516aa"
$ws.Range("L4").Value = "STD_040516aa"

# Move the current selection to G4 (was A6).
$ws.Range("G4").Select()
